$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$rng = $ws.Range("L2")
for ($i = 0; $i -lt 3; $i++) {
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Interior.ColorIndex = 9
}
